$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = "D2"; Value = "66.340.85" },
    @{ Cell = "E2"; Value = "  -3.83%  " },
    @{ Cell = "D3"; Value = "3.546.11" },
    @{ Cell = "E3"; Value = "  -4.64%  " },
    @{ Cell = "E4"; Value = "  +0.12%  " },
    @{ Cell = "D5"; Value = "575.81" },
    @{ Cell = "E5"; Value = "  -6.02%  " },
    @{ Cell = "D6"; Value = "188.05" },
    @{ Cell = "E6"; Value = "  -1.79%  " },
    @{ Cell = "D7"; Value = "3.540.69" },
    @{ Cell = "E7"; Value = "  -4.71%  " },
    @{ Cell = "D8"; Value = "0.610" },
    @{ Cell = "E8"; Value = "  -4.19%  " },
    @{ Cell = "E9"; Value = "  -0.13%  " },
    @{ Cell = "D10"; Value = "0.664" },
    @{ Cell = "E10"; Value = "  -7.30%  " },
    @{ Cell = "D11"; Value = "0.145" },
    @{ Cell = "E11"; Value = "  -9.75%  " },
    @{ Cell = "D12"; Value = "52.57" },
    @{ Cell = "E12"; Value = "  -7.42%  " },
    @{ Cell = "D13"; Value = "0.0000257" },
    @{ Cell = "E13"; Value = "  -11.74%  " },
    @{ Cell = "D14"; Value = "9.72" },
    @{ Cell = "E14"; Value = "  -8.03%  " },
    @{ Cell = "D15"; Value = "4.109.68" },
    @{ Cell = "E15"; Value = "  -4.54%  " },
    @{ Cell = "D16"; Value = "3.547.50" },
    @{ Cell = "E16"; Value = "  -4.54%  " },
    @{ Cell = "E17"; Value = "  -1.11%  " },
    @{ Cell = "D18"; Value = "18.19" },
    @{ Cell = "E18"; Value = "  -5.80%  " },
    @{ Cell = "D19"; Value = "66.129.90" },
    @{ Cell = "E19"; Value = "  -3.86%  " },
    @{ Cell = "D20"; Value = "12.05" },
    @{ Cell = "E20"; Value = "  -6.91%  " },
    @{ Cell = "E21"; Value = "  -7.60%  " },
    @{ Cell = "D22"; Value = "390.98" },
    @{ Cell = "E22"; Value = "  -5.02%  " },
    @{ Cell = "D23"; Value = "4.28" },
    @{ Cell = "E23"; Value = "  -6.82%  " },
    @{ Cell = "D24"; Value = "85.01" },
    @{ Cell = "E24"; Value = "  -5.01%  " },
    @{ Cell = "D25"; Value = "10.97" },
    @{ Cell = "E25"; Value = "  +0.33%  " },
    @{ Cell = "E26"; Value = "  -5.73%  " },
    @{ Cell = "D27"; Value = "12.28" },
    @{ Cell = "E27"; Value = "  -4.52%  " },
    @{ Cell = "D28"; Value = "6.04" },
    @{ Cell = "E28"; Value = "  -0.05%  " },
    @{ Cell = "D29"; Value = "3.49" },
    @{ Cell = "E29"; Value = "  -7.25%  " },
    @{ Cell = "D30"; Value = "8.81" },
    @{ Cell = "E30"; Value = "  -8.99%  " },
    @{ Cell = "D31"; Value = "30.80" },
    @{ Cell = "E31"; Value = "  -7.04%  " },
    @{ Cell = "D32"; Value = "7.14" },
    @{ Cell = "E32"; Value = "  -2.69%  " },
    @{ Cell = "D33"; Value = "627.85" },
    @{ Cell = "E33"; Value = "  -0.57%  " },
    @{ Cell = "D34"; Value = "12.12" },
    @{ Cell = "E34"; Value = "  -4.80%  " },
    @{ Cell = "D35"; Value = "63.42" },
    @{ Cell = "E35"; Value = "  -3.37%  " },
    @{ Cell = "E36"; Value = "  -7.94%  " },
    @{ Cell = "D37"; Value = "41.13" },
    @{ Cell = "E37"; Value = "  -9.70%  " },
    @{ Cell = "E38"; Value = "  +0.04%  " },
    @{ Cell = "D39"; Value = "0.393" },
    @{ Cell = "E39"; Value = "  -5.47%  " },
    @{ Cell = "D40"; Value = "0.0₃0756" },
    @{ Cell = "E40"; Value = "  -7.75%  " },
    @{ Cell = "D41"; Value = "0.998" },
    @{ Cell = "E41"; Value = "  -0.23%  " },
    @{ Cell = "D43"; Value = "2.965.96" },
    @{ Cell = "E43"; Value = "  +3.20%  " },
    @{ Cell = "D44"; Value = "2.79" },
    @{ Cell = "E44"; Value = "  -8.46%  " },
    @{ Cell = "D45"; Value = "2.47" },
    @{ Cell = "E45"; Value = "  -5.56%  " },
    @{ Cell = "D46"; Value = "0.0403" },
    @{ Cell = "E46"; Value = "  -9.06%  " },
    @{ Cell = "D47"; Value = "3.11" },
    @{ Cell = "E47"; Value = "  +1.62%  " },
    @{ Cell = "E48"; Value = "  -7.68%  " },
    @{ Cell = "D49"; Value = "138.07" },
    @{ Cell = "E49"; Value = "  -2.38%  " },
    @{ Cell = "D50"; Value = "8.40" },
    @{ Cell = "E50"; Value = "  -7.33%  " },
    @{ Cell = "D51"; Value = "2.71" },
    @{ Cell = "E51"; Value = "  -2.13%  " }
)

foreach ($chg in $changes) {
    $rng = $ws.Range($chg.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $chg.Value
    $rng.ClearFormats()
}
